$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$provinces = @(
  "Abra",
  "Agusan Del Norte",
  "Agusan Del Sur",
  "Aklan",
  "Albay",
  "Antique",
  "Apayao",
  "Aurora",
  "Basilan",
  "Bataan",
  "Batangas",
  "Biliran",
  "Bohol",
  "Bukidnon",
  "Bulacan",
  "Cagayan",
  "Camarines Norte",
  "Camarines Sur",
  "Camiguin",
  "Capiz",
  "Catanduanes",
  "Cavite",
  "Cebu",
  "Compostela Valley",
  "Davao Del Norte",
  "Davao Del Sur",
  "Davao Oriental",
  "Guimaras",
  "Ifugao",
  "Ilocos Norte",
  "Ilocos Sur",
  "Isabela",
  "Kalinga",
  "La Union",
  "Laguna",
  "Lanao Del Norte",
  "Leyte",
  "Marinduque",
  "Masbate",
  "Misamis Oriental",
  "Negros Occidental",
  "North Cotabato",
  "Northern Samar",
  "Nueva Ecija",
  "Nueva Vizcaya",
  "Occidental Mindoro",
  "Oriental Mindoro",
  "Palawan",
  "Pampanga",
  "Pangasinan",
  "Quezon",
  "Rizal",
  "Romblon",
  "Samar",
  "Sorsogon",
  "South Cotabato",
  "Southern Leyte",
  "Sultan Kudarat",
  "Surigao Del Norte",
  "Tarlac",
  "Zambales",
  "Zamboanga Del Sur",
  "Zamboanga Sibugay"
)

$bvals = @(
  "Low",
  "High",
  "Low",
  "Low",
  "Mid",
  "Low",
  "Low",
  "Mid",
  "Low",
  "High",
  "High",
  "High",
  "Mid",
  "Low",
  "High",
  "Mid",
  "Mid",
  "Mid",
  "Mid",
  "High",
  "High",
  "High",
  "High",
  "Low",
  "Mid",
  "High",
  "Low",
  "High",
  "Low",
  "Mid",
  "Mid",
  "Low",
  "Low",
  "High",
  "High",
  "Mid",
  "High",
  "Mid",
  "Low",
  "High",
  "Mid",
  "Low",
  "Low",
  "Low",
  "High",
  "High",
  "Mid",
  "Low",
  "High",
  "Mid",
  "Mid",
  "High",
  "Mid",
  "Low",
  "Mid",
  "High",
  "Low",
  "Low",
  "Mid",
  "Mid",
  "High",
  "Mid",
  "Low"
)

$cvals = @(
  "Mid",
  "High",
  "High",
  "High",
  "Low",
  "High",
  "Low",
  "Mid",
  "Mid",
  "Low",
  "Low",
  "High",
  "High",
  "Low",
  "Low",
  "Mid",
  "High",
  "Mid",
  "High",
  "Mid",
  "Mid",
  "Mid",
  "Mid",
  "Low",
  "Low",
  "Low",
  "High",
  "High",
  "Low",
  "Low",
  "Low",
  "Mid",
  "Mid",
  "Low",
  "Low",
  "Mid",
  "High",
  "High",
  "High",
  "High",
  "Mid",
  "High",
  "High",
  "Mid",
  "Mid",
  "Mid",
  "High",
  "High",
  "Mid",
  "Mid",
  "Low",
  "Low",
  "High",
  "Mid",
  "High",
  "Low",
  "Mid",
  "High",
  "High",
  "Mid",
  "Low",
  "Low",
  "Mid"
)

for ($i = 0; $i -lt $provinces.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $provinces[$i]
    $ws.Cells.Item($row, 2).Value = $bvals[$i]
    $ws.Cells.Item($row, 3).Value = $cvals[$i]
}

$ws.Rows("65:69").Delete()
